# Update recomputed NATMI ligand-receptor metrics (new TPM-based values)
# for Mdk-Sdc4 sheet. Values taken from the refreshed pipeline output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    # Row 2
    "G2" = 1.988074333333333
    "H2" = 5.964223
    "I2" = 0.01657769708907969
    "J2" = 0.01657769708907968
    "M2" = 0.5373756666666667
    "N2" = 1.612127
    "O2" = 0.007472820128982582
    "P2" = 0.007472820128982581
    "Q2" = 1.068342770257889
    "R2" = 9.615084932321
    "S2" = 0.0001238821484994507
    "T2" = 0.0001238821484994506
    # Row 3
    "G3" = 1.988074333333333
    "H3" = 5.964223
    "I3" = 0.01657769708907969
    "J3" = 0.01657769708907968
    "O3" = 0.1537223653287423
    "P3" = 0.1537223653287423
    "Q3" = 21.97673365493722
    "R3" = 197.790602894435
    "S3" = 0.002548362808236736
    "T3" = 0.002548362808236735
    # Row 4
    "G4" = 1.988074333333333
    "H4" = 5.964223
    "I4" = 0.01657769708907969
    "J4" = 0.01657769708907968
    "M4" = 30.561198
    "N4" = 91.683594
    "O4" = 0.4249882340167162
    "P4" = 0.4249882340167161
    "Q4" = 60.757933339718
    "R4" = 546.8214000574619
    "S4" = 0.007045326209952033
    "T4" = 0.007045326209952031
    # Row 5
    "G5" = 1.988074333333333
    "H5" = 5.964223
    "I5" = 0.01657769708907969
    "J5" = 0.01657769708907968
    "M5" = 29.75783666666667
    "N5" = 89.27351
    "O5" = 0.4138165805255589
    "P5" = 0.4138165805255589
    "Q5" = 59.16079129252555
    "R5" = 532.4471216327299
    "S5" = 0.006860125922391468
    "T5" = 0.006860125922391467
    # Row 6
    "I6" = 0.7746030815641455
    "J6" = 0.7746030815641454
    "M6" = 0.5373756666666667
    "N6" = 1.612127
    "O6" = 0.007472820128982582
    "P6" = 0.007472820128982581
    "Q6" = 49.91897231332978
    "R6" = 449.270750819968
    "S6" = 0.005788469499884484
    "T6" = 0.005788469499884482
    # Row 7
    "I7" = 0.7746030815641455
    "J7" = 0.7746030815641454
    "O7" = 0.1537223653287423
    "P7" = 0.1537223653287423
    "S7" = 0.1190738178889732
    "T7" = 0.1190738178889731
    # Row 8
    "I8" = 0.7746030815641455
    "J8" = 0.7746030815641454
    "M8" = 30.561198
    "N8" = 91.683594
    "O8" = 0.4249882340167162
    "P8" = 0.4249882340167161
    "Q8" = 2838.951764018944
    "R8" = 25550.56587617049
    "S8" = 0.3291971956978526
    "T8" = 0.3291971956978525
    # Row 9
    "I9" = 0.7746030815641455
    "J9" = 0.7746030815641454
    "M9" = 29.75783666666667
    "N9" = 89.27351
    "O9" = 0.4138165805255589
    "P9" = 0.4138165805255589
    "Q9" = 2764.324320604871
    "R9" = 24878.91888544384
    "S9" = 0.3205435984774353
    "T9" = 0.3205435984774352
    # Row 10
    "G10" = 23.741365
    "H10" = 71.22409500000001
    "I10" = 0.1979690350870239
    "J10" = 0.1979690350870239
    "M10" = 0.5373756666666667
    "N10" = 1.612127
    "O10" = 0.007472820128982582
    "P10" = 0.007472820128982581
    "Q10" = 12.75803184445167
    "R10" = 114.822286600065
    "S10" = 0.001479386990313571
    "T10" = 0.001479386990313571
    # Row 11
    "G11" = 23.741365
    "H11" = 71.22409500000001
    "I11" = 0.1979690350870239
    "J11" = 0.1979690350870239
    "O11" = 0.1537223653287423
    "P11" = 0.1537223653287423
    "Q11" = 262.4437358611417
    "R11" = 2361.993622750275
    "S11" = 0.0304322683354261
    "T11" = 0.0304322683354261
    # Row 12
    "G12" = 23.741365
    "H12" = 71.22409500000001
    "I12" = 0.1979690350870239
    "J12" = 0.1979690350870239
    "M12" = 30.561198
    "N12" = 91.683594
    "O12" = 0.4249882340167162
    "P12" = 0.4249882340167161
    "Q12" = 725.5645565552701
    "R12" = 6530.08100899743
    "S12" = 0.08413451061162763
    "T12" = 0.08413451061162762
    # Row 13
    "G13" = 23.741365
    "H13" = 71.22409500000001
    "I13" = 0.1979690350870239
    "J13" = 0.1979690350870239
    "M13" = 29.75783666666667
    "N13" = 89.27351
    "O13" = 0.4138165805255589
    "P13" = 0.4138165805255589
    "Q13" = 706.4916619137167
    "R13" = 6358.42495722345
    "S13" = 0.08192286914965664
    "T13" = 0.08192286914965664
    # Row 14
    "G14" = 1.301204666666667
    "H14" = 3.903614
    "I14" = 0.01085018625975097
    "J14" = 0.01085018625975097
    "M14" = 0.5373756666666667
    "N14" = 1.612127
    "O14" = 0.007472820128982582
    "P14" = 0.007472820128982581
    "Q14" = 0.6992357252197778
    "R14" = 6.293121526978001
    "S14" = 0.0000810814902850773
    "T14" = 0.00008108149028507727
    # Row 15
    "G15" = 1.301204666666667
    "H15" = 3.903614
    "I15" = 0.01085018625975097
    "J15" = 0.01085018625975097
    "O15" = 0.1537223653287423
    "P15" = 0.1537223653287423
    "Q15" = 14.38388289131445
    "R15" = 129.45494602183
    "S15" = 0.001667916296106339
    "T15" = 0.001667916296106339
    # Row 16
    "G16" = 1.301204666666667
    "H16" = 3.903614
    "I16" = 0.01085018625975097
    "J16" = 0.01085018625975097
    "M16" = 30.561198
    "N16" = 91.683594
    "O16" = 0.4249882340167162
    "P16" = 0.4249882340167161
    "Q16" = 39.76637345652401
    "R16" = 357.897361108716
    "S16" = 0.004611201497284004
    "T16" = 0.004611201497284004
    # Row 17
    "G17" = 1.301204666666667
    "H17" = 3.903614
    "I17" = 0.01085018625975097
    "J17" = 0.01085018625975097
    "M17" = 29.75783666666667
    "N17" = 89.27351
    "O17" = 0.4138165805255589
    "P17" = 0.4138165805255589
    "Q17" = 38.72103594057111
    "R17" = 348.48932346514
    "S17" = 0.004489986976075551
    "T17" = 0.004489986976075551
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
